# The workbook tracks one weekly price record per row for Espárragos
# (Feria Lagunitas de Puerto Montt). This commit adds one more weekly
# observation, which lands in the middle of the existing table (row 31),
# pushing the existing rows 31-63 down to 32-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31, shifting rows 31..63 down to 32..64
# (xlShiftDown = -4121).
$ws.Rows.Item(31).Insert(-4121)

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44894
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 300000000
$ws.Range("G31").Value = "Espárragos"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = 1800
$ws.Range("N31").Value = "$/kilo"
$ws.Range("O31").Value = "Provincia de Linares"
$ws.Range("P31").Value = 1800
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
